$d = $word.ActiveDocument

# Remove the table that held the old "{{ messages_dynamics }}" placeholder.
$table = $null
foreach ($t in $d.Tables) {
    if ($t.Range.Text -like "*messages_dynamics*") {
        $table = $t
        break
    }
}
if ($table -eq $null) {
    $table = $d.Tables.Item(1)
}
$table.Delete()

# The placeholder now belongs in the paragraph that used to be the empty
# trailing paragraph after the table. Re-read paragraphs fresh (post-delete)
# via the Content range so we get up to date Range positions.
$paragraphs = $d.Content.Paragraphs
$target = $paragraphs.Item($paragraphs.Count)

# Match the cell paragraph's old widow/orphan control setting.
$target.Format.WidowControl = $false

# Add the placeholder text into the (currently empty) run of that paragraph.
$target.Range.InsertAfter("{{ messages_dynamics }}")
